# Apply the "sharebrowser module" update to the arbeitspakete sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arbeitspakete")

# 1) Extend the existing TODO note in F7 with the new remark about the
#    smbbrowser module.
$ws.Range("F7").Value2 = "TODO: prüfe Client-Namen; übergebe User/Password von App an Subsysteme --> smbbrowser-modul"

# 2) Bump the estimated real effort for that row (9 -> 12 h) and let the
#    dependent SUM formula in row 19 recalc on its own.
$ws.Range("D7").Value2 = 12

# 3) Row 7 grew taller to fit the longer note.
$ws.Rows.Item(7).RowHeight = 46.25

# 4) Add the new work package row (33) describing the wizard follow-up task.
$ws.Range("A33").Value2 = "Wizard für user/password basierte Modulauswahl"
$ws.Range("B33").Value2 = 16
$ws.Range("C33").Value2 = 0.25
$ws.Range("D33").Value2 = 4

# Match number formatting (percentage) used by the other rows' "% erreicht"
# column, since a bare numeric assignment would otherwise inherit the plain
# column default style instead of the percentage style.
$ws.Range("C31").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null

# 5) Update the last-used selection to reflect where the author left off.
$ws.Range("C7").Select() | Out-Null
